# Apply updated Leve profit/price figures across the Ifrit_Profits sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 19
$ws.Range("H19").Value = 2201.95
$ws.Range("I19").Value = 3086.6365
$ws.Range("J19").Value = 1120.6666
$ws.Range("K19").Value = 3086.6365
$ws.Range("L19").Value = 1120.6666
$ws.Range("M19").Value = -2911.6365
$ws.Range("N19").Value = -1470.6666
# Row 40
$ws.Range("H40").Value = 1049.2222
$ws.Range("I40").Value = 994.8333
$ws.Range("J40").Value = 1158
$ws.Range("K40").Value = 994.8333
$ws.Range("L40").Value = 1158
$ws.Range("M40").Value = -819.8333
$ws.Range("N40").Value = -1508
# Row 64
$ws.Range("H64").Value = 4440
$ws.Range("I64").Value = 5000
$ws.Range("J64").Value = 4066.6667
$ws.Range("K64").Value = 5000
$ws.Range("L64").Value = 4066.6667
$ws.Range("M64").Value = -4752
$ws.Range("N64").Value = -4562.6667
# Row 67
$ws.Range("H67").Value = 4440
$ws.Range("I67").Value = 5000
$ws.Range("J67").Value = 4066.6667
$ws.Range("K67").Value = 5000
$ws.Range("L67").Value = 4066.6667
$ws.Range("M67").Value = -4142
$ws.Range("N67").Value = -5782.6667
# Row 106
$ws.Range("H106").Value = 2219
$ws.Range("I106").Value = 2257.0833
$ws.Range("J106").Value = 2066.6667
$ws.Range("K106").Value = 2257.0833
$ws.Range("L106").Value = 2066.6667
$ws.Range("M106").Value = -1626.0833
$ws.Range("N106").Value = -3328.6667
# Row 112
$ws.Range("H112").Value = 52632810
$ws.Range("J112").Value = 71429990
$ws.Range("L112").Value = 214289970
$ws.Range("N112").Value = -214292186
# Row 124
$ws.Range("H124").Value = 43000
$ws.Range("J124").Value = 43000
$ws.Range("L124").Value = 43000
$ws.Range("N124").Value = -52820
# Row 125
$ws.Range("H125").Value = 1156.25
$ws.Range("I125").Value = 1130
$ws.Range("J125").Value = 1200
$ws.Range("K125").Value = 10170
$ws.Range("L125").Value = 10800
$ws.Range("M125").Value = -7710
$ws.Range("N125").Value = -15720
# Row 137
$ws.Range("H137").Value = 1719.2858
$ws.Range("I137").Value = 1227.0588
$ws.Range("J137").Value = 2184.1667
$ws.Range("K137").Value = 3681.1764
$ws.Range("L137").Value = 6552.500100000001
$ws.Range("M137").Value = -1131.1764
$ws.Range("N137").Value = -11652.5001
# Row 139
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()
# Row 141
$ws.Range("H141").Value = 1819
$ws.Range("I141").Value = 1648.75
$ws.Range("K141").Value = 4946.25
$ws.Range("M141").Value = 233.75

$ws = $wb.Worksheets.Item("ARM")
# Row 92
$ws.Range("H92").Value = 23275
$ws.Range("J92").Value = 23275
$ws.Range("L92").Value = 23275
$ws.Range("N92").Value = -28267

$ws = $wb.Worksheets.Item("BSM")
# Row 9
$ws.Range("H9").Value = 14800
$ws.Range("J9").Value = 14800
$ws.Range("L9").Value = 14800
$ws.Range("N9").Value = -15136
# Row 59
$ws.Range("H59").Value = 28926.666
$ws.Range("J59").Value = 28926.666
$ws.Range("L59").Value = 28926.666
$ws.Range("N59").Value = -30620.666
# Row 94
$ws.Range("H94").Value = 1007.3333
$ws.Range("I94").Value = 1023.7778
$ws.Range("J94").Value = 859.3333
$ws.Range("K94").Value = 1023.7778
$ws.Range("L94").Value = 859.3333
$ws.Range("M94").Value = -572.7778
$ws.Range("N94").Value = -1761.3333
# Row 107
$ws.Range("H107").Value = 1544.2307
$ws.Range("I107").Value = 1508.6
$ws.Range("J107").Value = 1617.5883
$ws.Range("K107").Value = 1508.6
$ws.Range("L107").Value = 1617.5883
$ws.Range("M107").Value = 411.4000000000001
$ws.Range("N107").Value = -5457.588299999999

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 1339
$ws.Range("I31").Value = 1134.7273
$ws.Range("J31").Value = 2462.5
$ws.Range("K31").Value = 1134.7273
$ws.Range("L31").Value = 2462.5
$ws.Range("M31").Value = -839.7273
$ws.Range("N31").Value = -3052.5
# Row 34
$ws.Range("H34").Value = 1339
$ws.Range("I34").Value = 1134.7273
$ws.Range("J34").Value = 2462.5
$ws.Range("K34").Value = 1134.7273
$ws.Range("L34").Value = 2462.5
$ws.Range("M34").Value = -932.7273
$ws.Range("N34").Value = -2866.5
# Row 132
$ws.Range("H132").Value = 3168.238
$ws.Range("I132").Value = 2897.2144
$ws.Range("K132").Value = 8691.643199999999
$ws.Range("M132").Value = -6161.643199999999
# Row 134
$ws.Range("H134").Value = 4340.7334
$ws.Range("I134").Value = 4465.0713
$ws.Range("J134").Value = 2600
$ws.Range("K134").Value = 13395.2139
$ws.Range("L134").Value = 7800
$ws.Range("M134").Value = -10860.2139
$ws.Range("N134").Value = -12870

$ws = $wb.Worksheets.Item("CUL")
# Row 2
$ws.Range("H2").Value = 1424530.5
$ws.Range("I2").Value = 4.6666665
$ws.Range("J2").Value = 2136793.2
$ws.Range("K2").Value = 27.999999
$ws.Range("L2").Value = 12820759.2
$ws.Range("M2").Value = 85.000001
$ws.Range("N2").Value = -12820985.2
# Row 47
$ws.Range("H47").Value = 111.181816
$ws.Range("I47").Value = 112.3
$ws.Range("J47").Value = 100
$ws.Range("K47").Value = 336.9
$ws.Range("L47").Value = 300
$ws.Range("M47").Value = 94.10000000000002
$ws.Range("N47").Value = -1162
# Row 113
$ws.Range("H113").Value = 596.7826
$ws.Range("J113").Value = 586.0625
$ws.Range("L113").Value = 1758.1875
$ws.Range("N113").Value = -6098.1875
# Row 131
$ws.Range("H131").Value = 1756947.2
$ws.Range("I131").Value = 5913.3335
$ws.Range("J131").Value = 2085266
$ws.Range("K131").Value = 17740.0005
$ws.Range("L131").Value = 6255798
$ws.Range("M131").Value = -12700.0005
$ws.Range("N131").Value = -6265878
# Row 133
$ws.Range("H133").Value = 6575.2334
$ws.Range("I133").Value = 3358.75
$ws.Range("J133").Value = 7744.864
$ws.Range("K133").Value = 10076.25
$ws.Range("L133").Value = 23234.592
$ws.Range("M133").Value = -5016.25
$ws.Range("N133").Value = -33354.592

$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 101924.5
$ws.Range("I80").Value = 2499.5
$ws.Range("J80").Value = 121809.5
$ws.Range("K80").Value = 2499.5
$ws.Range("L80").Value = 121809.5
$ws.Range("M80").Value = -1501.5
$ws.Range("N80").Value = -123805.5
# Row 83
$ws.Range("H83").Value = 101924.5
$ws.Range("I83").Value = 2499.5
$ws.Range("J83").Value = 121809.5
$ws.Range("K83").Value = 12497.5
$ws.Range("L83").Value = 609047.5
$ws.Range("M83").Value = -7505.5
$ws.Range("N83").Value = -619031.5
# Row 107
$ws.Range("H107").Value = 737.6923
$ws.Range("I107").Value = 438
$ws.Range("J107").Value = 827.6
$ws.Range("K107").Value = 438
$ws.Range("L107").Value = 827.6
$ws.Range("M107").Value = 1482
$ws.Range("N107").Value = -4667.6
# Row 122
$ws.Range("H122").Value = 3179.6428
$ws.Range("I122").Value = 3167.9167
$ws.Range("J122").Value = 3250
$ws.Range("K122").Value = 9503.750100000001
$ws.Range("L122").Value = 9750
$ws.Range("M122").Value = -7053.750100000001
$ws.Range("N122").Value = -14650

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 2084.7144
$ws.Range("I7").Value = 1938.6
$ws.Range("J7").Value = 2450
$ws.Range("K7").Value = 1938.6
$ws.Range("L7").Value = 2450
$ws.Range("M7").Value = -1826.6
$ws.Range("N7").Value = -2674
# Row 16
$ws.Range("H16").Value = 1747.72
$ws.Range("I16").Value = 1816.375
$ws.Range("J16").Value = 100
$ws.Range("K16").Value = 1816.375
$ws.Range("L16").Value = 100
$ws.Range("M16").Value = -1646.375
$ws.Range("N16").Value = -440
# Row 22
$ws.Range("H22").Value = 323.76923
$ws.Range("I22").Value = 283.22223
$ws.Range("J22").Value = 415
$ws.Range("K22").Value = 283.22223
$ws.Range("L22").Value = 415
$ws.Range("M22").Value = 11.77776999999998
$ws.Range("N22").Value = -1005
# Row 27
$ws.Range("H27").Value = 323.76923
$ws.Range("I27").Value = 283.22223
$ws.Range("J27").Value = 415
$ws.Range("K27").Value = 283.22223
$ws.Range("L27").Value = 415
$ws.Range("M27").Value = -176.22223
$ws.Range("N27").Value = -629
# Row 40
$ws.Range("H40").Value = 1979.5625
$ws.Range("I40").Value = 1833.909
$ws.Range("J40").Value = 2300
$ws.Range("K40").Value = 1833.909
$ws.Range("L40").Value = 2300
$ws.Range("M40").Value = -1697.909
$ws.Range("N40").Value = -2572
# Row 114
$ws.Range("H114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("L114").Value = 0
$ws.Range("N114").ClearContents()
# Row 122
$ws.Range("H122").Value = 5627.757
$ws.Range("I122").Value = 6417.2856
$ws.Range("K122").Value = 19251.8568
$ws.Range("M122").Value = -16801.8568
# Row 126
$ws.Range("H126").Value = 2084.7144
$ws.Range("I126").Value = 1938.6
$ws.Range("J126").Value = 2450
$ws.Range("K126").Value = 5815.799999999999
$ws.Range("L126").Value = 7350
$ws.Range("M126").Value = -3345.799999999999
$ws.Range("N126").Value = -12290
# Row 136
$ws.Range("H136").Value = 2453.1667
$ws.Range("I136").Value = 1001.3333
$ws.Range("J136").Value = 3905
$ws.Range("K136").Value = 3003.9999
$ws.Range("L136").Value = 11715
$ws.Range("M136").Value = -453.9998999999998
$ws.Range("N136").Value = -16815

$ws = $wb.Worksheets.Item("WVR")
# Row 94
$ws.Range("H94").Value = 24000
$ws.Range("J94").Value = 24000
$ws.Range("L94").Value = 24000
$ws.Range("N94").Value = -25802
# Row 97
$ws.Range("H97").Value = 0
$ws.Range("J97").Value = 0
$ws.Range("L97").Value = 0
$ws.Range("N97").ClearContents()
